# Added RPA.Excel.Files test - Reading Writing data
#
# - Adds a "Result" header to the "login" sheet (new column C) and makes
#   that sheet the active tab with selection left on D2.
# - Sheet3 ("Data2") stops being the tab-selected sheet (handled
#   automatically once another sheet becomes active).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("login")
$ws.Activate()

$ws.Range("C1").Value = "Result"

$ws.Range("D2").Select() | Out-Null
